# rename Pmax to PmaxLink
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Header rename: F4 "pPmax" -> "pPmaxLink"
$ws.Range("F4").Value = "pPmaxLink"

# Update the Pmax data column (F8:F16) with the recalculated values and
# switch it to a one-decimal "0.0" number format (matching the precision
# now needed by the non-round numbers), keeping its existing fill color.
$ws.Range("F8").Value = 672.60519999999997
$ws.Range("F9").Value = 313.88240000000002
$ws.Range("F10").Value = 448.40350000000001
$ws.Range("F11").Value = 896.80690000000004
$ws.Range("F12").Value = 1661.2634
$ws.Range("F13").Value = 2028.3655000000001
$ws.Range("F14").Value = 358.72280000000001
$ws.Range("F15").Value = 2300.5331000000001
$ws.Range("F16").Value = 224.20169999999999

$ws.Range("F8:F16").NumberFormat = "0.0"

# Move the active selection to F4 (matches the new sheetView selection)
$ws.Range("F4").Select()
